$d = $word.ActiveDocument

# 1) Simplify "( Parts A & B)" -> "(Parts A & B)"
#    Collapses the 3 separate runs (space / "( Parts" / " A & ") into a single
#    run's text, which also removes the now-orphaned gramStart/gramEnd proofErr
#    markers around "( Parts".
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(" ( Parts A & ", $true, $false, $false, $false, $false, $true, 1, $false, " (Parts A & ", 2)

# 2) After the "MDN: .../Global_Objects/Array" paragraph, add a new empty
#    paragraph with the same indentation (ind left=426).
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Global_Objects/Array*") {
        $p.Range.InsertParagraphAfter()
        break
    }
}
